# Auto-generated COM-interop script reproducing the AggTagTemplate.xlsx diff.
# Adds 5 new worksheets (Msd, Rollup, Rollups, Cube, GroupingSets) with
# header/data rows styled to match the JETT 'agg' tag demo template.

$wb = $excel.ActiveWorkbook

# Blue fill color for header rows: RGB 99CCFF, encoded BGR for the COM Color property.
$headerFillColor = 16764057
$currencyFormat = '"$"#,##0.00'

### Sheet: Msd ###
$lastIdx = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = 'Msd'

$ws.Columns.Item(1).ColumnWidth = 25.6
$ws.Columns.Item(2).ColumnWidth = 11.45

$headerRange = $ws.Range('A1:B1')
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = $headerFillColor
$ws.Range('A1').Value = 'Is A Manager'
$ws.Range('B1').Value = 'Total Salary'

$ws.Range('A2').Borders.LineStyle = 1
$lastCell = $ws.Range('B2')
$lastCell.Borders.LineStyle = 1
$lastCell.NumberFormat = $currencyFormat

$ws.Range('A2').Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager()" useMsd="true"><jt:forEach items="${values}" var="value">${value.getPropertyValue(0)}'
$ws.Range('B2').Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'

$ws.PageSetup.Orientation = 1


### Sheet: Rollup ###
$lastIdx = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = 'Rollup'

$ws.Columns.Item(1).ColumnWidth = 25.6
$ws.Columns.Item(2).ColumnWidth = 25.6
$ws.Columns.Item(3).ColumnWidth = 11.45

$headerRange = $ws.Range('A1:C1')
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = $headerFillColor
$ws.Range('A1').Value = 'Is A Manager'
$ws.Range('B1').Value = 'Title'
$ws.Range('C1').Value = 'Total Salary'

$ws.Range('A2').Borders.LineStyle = 1
$ws.Range('B2:B2').Borders.LineStyle = 1
$lastCell = $ws.Range('C2')
$lastCell.Borders.LineStyle = 1
$lastCell.NumberFormat = $currencyFormat

$ws.Range('A2').Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title" rollup="${[0, 1]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
$ws.Range('B2').Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
$ws.Range('C2').Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'


### Sheet: Rollups ###
$lastIdx = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = 'Rollups'

$ws.Columns.Item(1).ColumnWidth = 25.6
$ws.Columns.Item(2).ColumnWidth = 25.6
$ws.Columns.Item(3).ColumnWidth = 30.15
$ws.Columns.Item(4).ColumnWidth = 11.45

$headerRange = $ws.Range('A1:D1')
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = $headerFillColor
$ws.Range('A1').Value = 'Is A Manager'
$ws.Range('B1').Value = 'Title'
$ws.Range('C1').Value = 'Catch Phrase'
$ws.Range('D1').Value = 'Total Salary'

$ws.Range('A2').Borders.LineStyle = 1
$ws.Range('B2:C2').Borders.LineStyle = 1
$lastCell = $ws.Range('D2')
$lastCell.Borders.LineStyle = 1
$lastCell.NumberFormat = $currencyFormat

$ws.Range('A2').Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title;catchPhrase" rollups="${[[1], [2]]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1);getPropertyValue(2)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
$ws.Range('B2').Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
$ws.Range('C2').Value = '${value.isGrouping(2) ? ''All Values'' : value.getPropertyValue(2)}'
$ws.Range('D2').Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'


### Sheet: Cube ###
$lastIdx = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = 'Cube'

$ws.Columns.Item(1).ColumnWidth = 25.6
$ws.Columns.Item(2).ColumnWidth = 25.6
$ws.Columns.Item(3).ColumnWidth = 30.15
$ws.Columns.Item(4).ColumnWidth = 11.45

$headerRange = $ws.Range('A1:D1')
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = $headerFillColor
$ws.Range('A1').Value = 'Is A Manager'
$ws.Range('B1').Value = 'Title'
$ws.Range('C1').Value = 'Catch Phrase'
$ws.Range('D1').Value = 'Total Salary'

$ws.Range('A2').Borders.LineStyle = 1
$ws.Range('B2:C2').Borders.LineStyle = 1
$lastCell = $ws.Range('D2')
$lastCell.Borders.LineStyle = 1
$lastCell.NumberFormat = $currencyFormat

$ws.Range('A2').Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title;catchPhrase" cube="${[0, 1, 2]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1);getPropertyValue(2)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
$ws.Range('B2').Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
$ws.Range('C2').Value = '${value.isGrouping(2) ? ''All Values'' : value.getPropertyValue(2)}'
$ws.Range('D2').Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'


### Sheet: GroupingSets ###
$lastIdx = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = 'GroupingSets'

$ws.Columns.Item(1).ColumnWidth = 25.6
$ws.Columns.Item(2).ColumnWidth = 25.6
$ws.Columns.Item(3).ColumnWidth = 30.15
$ws.Columns.Item(4).ColumnWidth = 11.45

$headerRange = $ws.Range('A1:D1')
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = $headerFillColor
$ws.Range('A1').Value = 'Is A Manager'
$ws.Range('B1').Value = 'Title'
$ws.Range('C1').Value = 'Catch Phrase'
$ws.Range('D1').Value = 'Total Salary'

$ws.Range('A2').Borders.LineStyle = 1
$ws.Range('B2:C2').Borders.LineStyle = 1
$lastCell = $ws.Range('D2')
$lastCell.Borders.LineStyle = 1
$lastCell.NumberFormat = $currencyFormat

$ws.Range('A2').Value = '<jt:agg items="${employees}" aggs="Sum(salary)" valuesVar="values" groupBy="isManager();title;catchPhrase" groupingSets="${[[0], [1, 2]]}"><jt:forEach items="${values}" var="value" orderBy="getPropertyValue(0);getPropertyValue(1);getPropertyValue(2)">${value.isGrouping(0) ? ''All Values'' : value.getPropertyValue(0)}'
$ws.Range('B2').Value = '${value.isGrouping(1) ? ''All Values'' : value.getPropertyValue(1)}'
$ws.Range('C2').Value = '${value.isGrouping(2) ? ''All Values'' : value.getPropertyValue(2)}'
$ws.Range('D2').Value = '${value.getAggregateValue(0)}</jt:forEach></jt:agg>'


Write-Host "Workbook now has $($wb.Worksheets.Count) sheets"
